# Update "想去人数" (want-to-go count) figures to the latest scraped values.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F9").Value = 145
$wsExpo.Range("F10").Value = 2376
$wsExpo.Range("F11").Value = 14

# Sheet "演出" (Performances)
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F3").Value = 44

# Sheet "全部类型" (All types / combined)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F10").Value = 145
$wsAll.Range("F11").Value = 2377
$wsAll.Range("F12").Value = 14
$wsAll.Range("F13").Value = 44
